$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1504.0714
$ws.Range("I129").Value = 719.0833
$ws.Range("K129").Value = 2157.2499
$ws.Range("M129").Value = 2842.7501
$ws.Range("H132").Value = 106720.625
$ws.Range("I132").Value = 229991.83
$ws.Range("J132").Value = 14267.23
$ws.Range("K132").Value = 689975.49
$ws.Range("L132").Value = 42801.69
$ws.Range("M132").Value = -687445.49
$ws.Range("N132").Value = -47861.69
$ws.Range("H138").Value = 6137.75
$ws.Range("I138").Value = 2720.7778
$ws.Range("J138").Value = 7016.4
$ws.Range("K138").Value = 8162.3334
$ws.Range("L138").Value = 21049.2
$ws.Range("M138").Value = -3022.3334
$ws.Range("N138").Value = -31329.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4754.793
$ws.Range("I32").Value = 4111.4585
$ws.Range("J32").Value = 7842.8
$ws.Range("K32").Value = 4111.4585
$ws.Range("L32").Value = 7842.8
$ws.Range("M32").Value = -3824.4585
$ws.Range("N32").Value = -8416.799999999999
$ws.Range("H43").Value = 8455.200000000001
$ws.Range("J43").Value = 8944
$ws.Range("L43").Value = 8944
$ws.Range("N43").Value = -9570
$ws.Range("H61").Value = 11642.405
$ws.Range("J61").Value = 10638.667
$ws.Range("L61").Value = 10638.667
$ws.Range("N61").Value = -11062.667
$ws.Range("H74").Value = 5437942.5
$ws.Range("I74").Value = 11364855
$ws.Range("K74").Value = 11364855
$ws.Range("M74").Value = -11363981
$ws.Range("H77").Value = 5437942.5
$ws.Range("I77").Value = 11364855
$ws.Range("K77").Value = 56824275
$ws.Range("M77").Value = -56819907
$ws.Range("H102").Value = 572054.5600000001
$ws.Range("I102").Value = 572054.5600000001
$ws.Range("K102").Value = 572054.5600000001
$ws.Range("M102").Value = -570432.5600000001
$ws.Range("H122").Value = 4047.7307
$ws.Range("I122").Value = 2060.8235
$ws.Range("K122").Value = 6182.470499999999
$ws.Range("M122").Value = -3732.470499999999
$ws.Range("H132").Value = 4249.9033
$ws.Range("I132").Value = 3525.2163
$ws.Range("J132").Value = 5322.44
$ws.Range("K132").Value = 10575.6489
$ws.Range("L132").Value = 15967.32
$ws.Range("M132").Value = -8045.6489
$ws.Range("N132").Value = -21027.32
$ws.Range("H136").Value = 11642.405
$ws.Range("J136").Value = 10638.667
$ws.Range("L136").Value = 31916.001
$ws.Range("N136").Value = -37016.001
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 69500
$ws.Range("J92").Value = 69500
$ws.Range("L92").Value = 69500
$ws.Range("N92").Value = -74492
$ws.Range("H99").Value = 34032150
$ws.Range("I99").Value = 170140770
$ws.Range("K99").Value = 170140770
$ws.Range("M99").Value = -170139272
$ws.Range("H105").Value = 150001620
$ws.Range("I105").Value = 166668260
$ws.Range("J105").Value = 1795
$ws.Range("K105").Value = 166668260
$ws.Range("L105").Value = 1795
$ws.Range("M105").Value = -166666513
$ws.Range("N105").Value = -5289
$ws.Range("H107").Value = 1556
$ws.Range("I107").Value = 1213.7059
$ws.Range("J107").Value = 2719.8
$ws.Range("K107").Value = 1213.7059
$ws.Range("L107").Value = 2719.8
$ws.Range("M107").Value = 706.2941000000001
$ws.Range("N107").Value = -6559.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29413854
$ws.Range("I31").Value = 35715304
$ws.Range("K31").Value = 35715304
$ws.Range("M31").Value = -35715009
$ws.Range("H34").Value = 29413854
$ws.Range("I34").Value = 35715304
$ws.Range("K34").Value = 35715304
$ws.Range("M34").Value = -35715102
$ws.Range("H58").Value = 627457.3
$ws.Range("I58").Value = 835076.8
$ws.Range("K58").Value = 835076.8
$ws.Range("M58").Value = -834873.8
$ws.Range("H132").Value = 49392776
$ws.Range("I132").Value = 57973780
$ws.Range("K132").Value = 173921340
$ws.Range("M132").Value = -173918810
$ws.Range("H134").Value = 1718.9
$ws.Range("I134").Value = 1524.375
$ws.Range("J134").Value = 2497
$ws.Range("K134").Value = 4573.125
$ws.Range("L134").Value = 7491
$ws.Range("M134").Value = -2038.125
$ws.Range("N134").Value = -12561
$ws.Range("H136").Value = 627457.3
$ws.Range("I136").Value = 835076.8
$ws.Range("K136").Value = 2505230.4
$ws.Range("M136").Value = -2502680.4
$ws.Range("H141").Value = 127400.336
$ws.Range("J141").Value = 127400.336
$ws.Range("L141").Value = 127400.336
$ws.Range("N141").Value = -137760.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5738496
$ws.Range("I4").Value = 2600095
$ws.Range("K4").Value = 7800285
$ws.Range("M4").Value = -7800173
$ws.Range("H37").Value = 83466240
$ws.Range("J37").Value = 83466240
$ws.Range("L37").Value = 250398720
$ws.Range("N37").Value = -250398944
$ws.Range("H86").Value = 340
$ws.Range("I86").Value = 197.5
$ws.Range("J86").Value = 625
$ws.Range("K86").Value = 592.5
$ws.Range("L86").Value = 1875
$ws.Range("M86").Value = 593.5
$ws.Range("N86").Value = -4247
$ws.Range("H89").Value = 340
$ws.Range("I89").Value = 197.5
$ws.Range("J89").Value = 625
$ws.Range("K89").Value = 1777.5
$ws.Range("L89").Value = 5625
$ws.Range("M89").Value = 4150.5
$ws.Range("N89").Value = -17481
$ws.Range("H128").Value = 339023.25
$ws.Range("I128").Value = 339023.25
$ws.Range("K128").Value = 1017069.75
$ws.Range("M128").Value = -1012089.75
$ws.Range("H131").Value = 16047863
$ws.Range("J131").Value = 7581504
$ws.Range("L131").Value = 22744512
$ws.Range("N131").Value = -22754592
$ws.Range("H140").Value = 6490.48
$ws.Range("J140").Value = 8851.799999999999
$ws.Range("L140").Value = 26555.4
$ws.Range("N140").Value = -36915.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15879345
$ws.Range("I70").Value = 15879345
$ws.Range("K70").Value = 15879345
$ws.Range("M70").Value = -15879075
$ws.Range("H73").Value = 15879345
$ws.Range("I73").Value = 15879345
$ws.Range("K73").Value = 15879345
$ws.Range("M73").Value = -15878409
$ws.Range("H122").Value = 504217.22
$ws.Range("I122").Value = 690242
$ws.Range("J122").Value = 8151.1665
$ws.Range("K122").Value = 2070726
$ws.Range("L122").Value = 24453.4995
$ws.Range("M122").Value = -2068276
$ws.Range("N122").Value = -29353.4995
$ws.Range("H126").Value = 3484.6
$ws.Range("I126").Value = 2193.2307
$ws.Range("J126").Value = 11878.5
$ws.Range("K126").Value = 6579.6921
$ws.Range("L126").Value = 35635.5
$ws.Range("M126").Value = -4109.6921
$ws.Range("N126").Value = -40575.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3574.6843
$ws.Range("I7").Value = 3359.7878
$ws.Range("K7").Value = 3359.7878
$ws.Range("M7").Value = -3247.7878
$ws.Range("H40").Value = 33337662
$ws.Range("I40").Value = 9095296
$ws.Range("J40").Value = 166670670
$ws.Range("K40").Value = 9095296
$ws.Range("L40").Value = 166670670
$ws.Range("M40").Value = -9095160
$ws.Range("N40").Value = -166670942
$ws.Range("H122").Value = 66670484
$ws.Range("I122").Value = 90912840
$ws.Range("K122").Value = 272738520
$ws.Range("M122").Value = -272736070
$ws.Range("H126").Value = 3574.6843
$ws.Range("I126").Value = 3359.7878
$ws.Range("K126").Value = 10079.3634
$ws.Range("M126").Value = -7609.3634
$ws.Range("H136").Value = 2039.8485
$ws.Range("I136").Value = 1835.5435
$ws.Range("J136").Value = 4725
$ws.Range("K136").Value = 5506.6305
$ws.Range("L136").Value = 14175
$ws.Range("M136").Value = -2956.6305
$ws.Range("N136").Value = -19275
$ws.Range("H139").Value = 90223.75
$ws.Range("J139").Value = 103632
$ws.Range("L139").Value = 103632
$ws.Range("N139").Value = -113912
$ws.Range("H140").Value = 95635.5
$ws.Range("J140").Value = 95635.5
$ws.Range("L140").Value = 95635.5
$ws.Range("N140").Value = -105995.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2086992.9
$ws.Range("I81").Value = 2318325.5
$ws.Range("K81").Value = 4636651
$ws.Range("M81").Value = -4635590
$ws.Range("H84").Value = 2086992.9
$ws.Range("I84").Value = 2318325.5
$ws.Range("K84").Value = 23183255
$ws.Range("M84").Value = -23177951
$ws.Range("H100").Value = 953727.9
$ws.Range("I100").Value = 1177261.9
$ws.Range("K100").Value = 2354523.8
$ws.Range("M100").Value = -2353982.8
$ws.Range("H122").Value = 2666.121
$ws.Range("I122").Value = 2700.2307
$ws.Range("J122").Value = 2539.4285
$ws.Range("K122").Value = 8100.6921
$ws.Range("L122").Value = 7618.2855
$ws.Range("M122").Value = -5650.6921
$ws.Range("N122").Value = -12518.2855
$ws.Range("H126").Value = 1352.4667
$ws.Range("I126").Value = 1164.0714
$ws.Range("K126").Value = 3492.2142
$ws.Range("M126").Value = -1022.2142
$ws.Range("H132").Value = 12082180
$ws.Range("I132").Value = 1393709.5
$ws.Range("K132").Value = 4181128.5
$ws.Range("M132").Value = -4178598.5
